$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.178.03"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "2.093.40"
$ws.Range("E3").Value = "  +2.79%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "229.34"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("E6").Value = "  +1.41%  "

$ws.Range("D7").Value = "60.85"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "0.0843"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").Value = "2.402.57"
$ws.Range("E12").Value = "  +2.65%  "

$ws.Range("D13").Value = "22.40"
$ws.Range("E13").Value = "  +5.05%  "

$ws.Range("D14").Value = "14.65"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").Value = "5.53"
$ws.Range("E15").Value = "  +6.87%  "

$ws.Range("D16").Value = "0.775"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").Value = "2.095.40"
$ws.Range("E17").Value = "  +3.16%  "

$ws.Range("D18").Value = "38.113.49"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "6.01"
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "70.28"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").Value = "224.43"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +2.99%  "

$ws.Range("D26").Value = "170.20"
$ws.Range("E26").Value = "  +2.23%  "

$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("D29").Value = "18.99"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("E30").Value = "  +5.34%  "

$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "2.39"
$ws.Range("E32").Value = "  +7.36%  "

$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").Value = "0.0605"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").Value = "6.51"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("E37").Value = "  +5.23%  "

$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  +7.92%  "

$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").Value = "18.07"
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("D41").Value = "1.559.24"
$ws.Range("E41").Value = "  +1.56%  "

$ws.Range("D42").Value = "100.32"
$ws.Range("E42").Value = "  +4.16%  "

$ws.Range("D43").Value = "0.0220"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D45").Value = "0.0906"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "4.13"
$ws.Range("E46").Value = "  +3.20%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("E49").Value = "  +2.46%  "

$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("D51").Value = "2.290.56"
$ws.Range("E51").Value = "  +2.75%  "
